# Update "Förändrad" (changed) date from 2023-09-11 (45180) to 2023-09-12 (45181)
# for every existing data row (rows 2-295, column C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C295").Value = 45181

# Row 295 gains an explicit row height (15pt, custom height) in the target file.
$ws.Rows.Item(295).RowHeight = 15

# Append the new row 296 with a fresh notification record.
# Start by cloning row 295's formatting down to row 296 so date columns keep
# their date number format and column R keeps its wrap-text style.
$ws.Range("A295:R295").Copy()
$ws.Range("A296:R296").PasteSpecial(-4122)

$ws.Range("A296").Value = "A 42308-2023"
$ws.Range("B296").Value = 45180
$ws.Range("C296").Value = 45181
$ws.Range("D296").Value = "DALARNAS LÄN"
$ws.Range("E296").Value = "VANSBRO"
$ws.Range("F296").Value = "Bergvik skog öst AB"
$ws.Range("G296").Value = 1.9
$ws.Range("H296").Value = 0
$ws.Range("I296").Value = 0
$ws.Range("J296").Value = 0
$ws.Range("K296").Value = 0
$ws.Range("L296").Value = 0
$ws.Range("M296").Value = 0
$ws.Range("N296").Value = 0
$ws.Range("O296").Value = 0
$ws.Range("P296").Value = 0
$ws.Range("Q296").Value = 0
